{"js": "// Clarify the text about splitting QPIs so the worked example matches current\n// practice: describe editing the QPI \"name\" column using the Years1:6 style\n// suffix (instead of literal year range), and note that the years (not just\n// \"the different\") are what make the data non-comparable.\n\nasync function replaceOnce(context, searchText, newText) {\n  const results = context.document.body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${searchText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"make the old QPI as QPI\" -> \"edit the old QPI name (ie 'QPI' columns) to\n//    be QPI\", and the worked example year range \"QPI 7: Nephron Sparing\n//    Surgery 2012:2017\" becomes the quoted \"Years1:6\" style name with a\n//    trailing explanation of what years that maps to.\nawait replaceOnce(\n  context,\n  \"then make the old QPI as QPI plus relevant years e.g., QPI 7: Nephron Sparing Surgery 2012:2017. Other\",\n  \"then edit the old QPI name (ie \\u2018QPI\\u2019 columns) to be QPI plus relevant years e.g., \\u201cQPI 7: Nephron Sparing Surgery: Years1:6\\u201d (for the years 2012 to 2017). Other\"\n);\n\n// 2) Wrap the archived-name worked example in quotation marks.\nawait replaceOnce(\n  context,\n  \"adding archived at end e.g. ZQPI 6: Neo-adjuvant Radiotherapy (archived).\",\n  \"adding archived at end e.g. \\u201cZQPI 6: Neo-adjuvant Radiotherapy (archived)\\u201d.\"\n);\n\n// 3) Clarify that it's the different *years* of data that aren't comparable.\nawait replaceOnce(\n  context,\n  \"Tableau will treat the data from the different as not being comparable\",\n  \"Tableau will treat the data from the different years as not being comparable\"\n);\n", "ps1": "# Clarify the text about splitting QPIs so the worked example matches current\n# practice: describe editing the QPI \"name\" column using the Years1:6 style\n# suffix (instead of literal year range), and note that the years (not just\n# \"the different\") are what make the data non-comparable.\n\n$d = $word.ActiveDocument\n\nfunction Replace-OnceInDocument($FindText, $ReplaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceText\n    $find.Forward = $true\n    $find.Wrap = 0          # wdFindStop\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # wdReplaceAll = 2\n    $result = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n    if (-not $result) {\n        throw \"Replace-OnceInDocument: could not find text: $FindText\"\n    }\n}\n\n# 1) \"make the old QPI as QPI\" -> \"edit the old QPI name (ie 'QPI' columns) to\n#    be QPI\", and the worked example year range \"QPI 7: Nephron Sparing\n#    Surgery 2012:2017\" becomes the quoted \"Years1:6\" style name with a\n#    trailing explanation of what years that maps to.\nReplace-OnceInDocument `\n    \"then make the old QPI as QPI plus relevant years e.g., QPI 7: Nephron Sparing Surgery 2012:2017. Other\" `\n    \"then edit the old QPI name (ie \u2018QPI\u2019 columns) to be QPI plus relevant years e.g., \u201cQPI 7: Nephron Sparing Surgery: Years1:6\u201d (for the years 2012 to 2017). Other\"\n\n# 2) Wrap the archived-name worked example in quotation marks.\nReplace-OnceInDocument `\n    \"adding archived at end e.g. ZQPI 6: Neo-adjuvant Radiotherapy (archived).\" `\n    \"adding archived at end e.g. \u201cZQPI 6: Neo-adjuvant Radiotherapy (archived)\u201d.\"\n\n# 3) Clarify that it's the different *years* of data that aren't comparable.\nReplace-OnceInDocument `\n    \"Tableau will treat the data from the different as not being comparable\" `\n    \"Tableau will treat the data from the different years as not being comparable\"\n"}
